$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'246.16"
$ws.Range("G2").Value = "'3"

# Row 3
$ws.Range("D3").Value = "'21.54"
$ws.Range("G3").Value = "'3"

# Row 4
$ws.Range("D4").Value = "'5.287"
$ws.Range("G4").Value = "'3"

# Row 5
$ws.Range("G5").Value = "'3"

# Row 6
$ws.Range("D6").Value = "'3.391"
$ws.Range("G6").Value = "'3"

# Row 7
$ws.Range("D7").Value = "'6.384"
$ws.Range("G7").Value = "'3"

# Row 8
$ws.Range("D8").Value = "'0.8177"
$ws.Range("G8").Value = "'3"

# Row 9
$ws.Range("D9").Value = "'0.9789"
$ws.Range("G9").Value = "'3"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1404"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "'3"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07432"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "'3"

# Row 12
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03152"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "'3"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03042"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "'3"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09293"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "'3"

# Row 15
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.561"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "'3"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001619"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "'3"

# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04719"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "'3"

# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005762"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "'3"

# Row 19
$ws.Range("D19").Value = "'0.006347"
$ws.Range("G19").Value = "'3"

# Row 20
$ws.Range("D20").Value = "'0.005054"
$ws.Range("G20").Value = "'3"

# Row 21
$ws.Range("G21").Value = "'3"

# Row 22
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("G22").Value = "'3"

# Row 23
$ws.Range("D23").Value = "'3.771"
$ws.Range("G23").Value = "'3"

# Row 24
$ws.Range("D24").Value = "'2.122"
$ws.Range("G24").Value = "'3"

# Row 25
$ws.Range("D25").Value = "'0.3252"
$ws.Range("G25").Value = "'3"

# Row 26
$ws.Range("D26").Value = "'0.1285"
$ws.Range("G26").Value = "'3"

# Row 27
$ws.Range("G27").Value = "'3"

# Row 28
$ws.Range("D28").Value = "'0.0003095"
$ws.Range("G28").Value = "'3"

# Row 29
$ws.Range("G29").Value = "'3"

# Row 30
$ws.Range("G30").Value = "'3"

# Row 31
$ws.Range("G31").Value = "'3"

# Row 32
$ws.Range("G32").Value = "'3"

# Row 33
$ws.Range("G33").Value = "'3"

# Row 34
$ws.Range("G34").Value = "'3"

# Row 35
$ws.Range("G35").Value = "'3"

# Row 36
$ws.Range("G36").Value = "'3"

# Row 37
$ws.Range("G37").Value = "'3"

# Row 38
$ws.Range("G38").Value = "'3"

# Row 39
$ws.Range("G39").Value = "'3"

# Row 40
$ws.Range("D40").Value = "'0.03925"
$ws.Range("G40").Value = "'3"

# Row 41
$ws.Range("D41").Value = "'0.007043"
$ws.Range("G41").Value = "'3"

# Row 42
$ws.Range("G42").Value = "'3"

# Row 43
$ws.Range("D43").Value = "'0.002848"
$ws.Range("G43").Value = "'3"

# Row 44
$ws.Range("D44").Value = "'0.007773"
$ws.Range("G44").Value = "'3"

# Row 45
$ws.Range("D45").Value = "'0.00005801"
$ws.Range("G45").Value = "'3"

# Row 46
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("G46").Value = "'3"

# Row 47
$ws.Range("D47").Value = "'0.0005492"
$ws.Range("G47").Value = "'3"

# Row 48
$ws.Range("D48").Value = "'0.6789"
$ws.Range("G48").Value = "'3"

# Row 49
$ws.Range("D49").Value = "'0.1430"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("G49").Value = "'3"

# Row 50
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("G50").Value = "'3"

# Row 51
$ws.Range("D51").Value = "'0.01008"
$ws.Range("G51").Value = "'3"
